# DATA TRANSPORT DAN POTONGAN.xlsx - apply "sweetaler 2 / tmt / tunjangan pangan" edit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns F..W alternate between the "style 5" (border-only) and "style 6"
# (font+border) cell formats that already exist in the sheet, starting with
# F = style5, G = style6, H = style5, ... W = style6.
# F and G already carry the correct style in every data row, so we only need
# to copy formats onto the columns whose style actually changes (H,J,L,N,P,R,T,V
# go from style6 -> style5); I,K,M,O,Q,S,U,W stay style6 (no-op, but harmless
# to leave alone).
$styleSourceOdd  = $ws.Range("F3")   # style 5 template (border only)
$styleSourceEven = $ws.Range("G3")   # style 6 template (font + border)

$columns = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W")

for ($r = 3; $r -le 7; $r++) {
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $col = $columns[$i]
        $addr = "$col$r"

        # Apply the alternating style by pasting formats from the matching
        # template cell in this same row (keeps borders/fonts/number formats
        # correct without fabricating brand-new style records).
        if ($i % 2 -eq 0) {
            $ws.Range("$col$r").Value = 0
            $styleSourceOdd.Copy()
            $ws.Range($addr).PasteSpecial(-4122)
        } else {
            $ws.Range("$col$r").Value = 0
            $styleSourceEven.Copy()
            $ws.Range($addr).PasteSpecial(-4122)
        }
    }
}
$excel.CutCopyMode = 0

# New row 11 / C11: an empty, quote-prefixed cell (as if a bare leading
# apostrophe were typed then the content cleared, leaving the text-quote
# formatting behind).
$ws.Range("C11").Value = "'"
$ws.Range("C11").ClearContents()

# Sheet view: drop the old horizontal scroll (topLeftCell="J1") and move the
# active selection to T19.
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
$ws.Range("T19").Select() | Out-Null
